$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.38%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.45%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.250"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.74%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07531"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'12.18%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.827"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.95%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.750"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'9.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.475"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7.09%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9135"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.49%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.01650"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2,451.60%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1696"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.08%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07548"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'6.85%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08040"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.33%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.02994"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.42%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09892"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'10.11%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001498"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.47%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04550"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.17%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006266"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.28%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.492"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.27%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.232"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.04%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'0.1342"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.471"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'14.40%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'4.56%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.91%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004443"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'1.75%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'19.41%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001737"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'7.41%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04504"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.29%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007255"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'7.18%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'8.95%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'0.69%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.63%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006211"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'8.13%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.870"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.72%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01297"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-13.56%"
$ws.Range("E47").Style = "Normal"
